# Auto-generated edit script applying the Halicarnassus_Profits.xlsx diff
# to the corresponding worksheets (ALC=sheet1 ... WVR=sheet8).
$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 9444.444
$ws.Range("I21").Value = 9444.444
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 9444.444
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -8976.444
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 9444.444
$ws.Range("I23").Value = 9444.444
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 9444.444
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -9210.444
$ws.Range("N23").ClearContents()
$ws.Range("H29").Value = 42.5
$ws.Range("I29").Value = 42.5
$ws.Range("K29").Value = 127.5
$ws.Range("M29").Value = 153.5
$ws.Range("H32").Value = 870
$ws.Range("J32").Value = 870
$ws.Range("L32").Value = 870
$ws.Range("N32").Value = -1522
$ws.Range("I43").Value = 2506.077
$ws.Range("J43").Value = 2434.6365
$ws.Range("K43").Value = 2506.077
$ws.Range("L43").Value = 2434.6365
$ws.Range("M43").Value = -2437.077
$ws.Range("N43").Value = -2572.6365
$ws.Range("H69").Value = 6982.4165
$ws.Range("I69").Value = 6800
$ws.Range("K69").Value = 20400
$ws.Range("M69").Value = -19526
$ws.Range("H72").Value = 6982.4165
$ws.Range("I72").Value = 6800
$ws.Range("K72").Value = 61200
$ws.Range("M72").Value = -56832
$ws.Range("H88").Value = 1226.6666
$ws.Range("I88").Value = 1258.125
$ws.Range("J88").Value = 975
$ws.Range("K88").Value = 1258.125
$ws.Range("L88").Value = 975
$ws.Range("M88").Value = -852.125
$ws.Range("N88").Value = -1787
$ws.Range("H91").Value = 1226.6666
$ws.Range("I91").Value = 1258.125
$ws.Range("J91").Value = 975
$ws.Range("K91").Value = 1258.125
$ws.Range("L91").Value = 975
$ws.Range("M91").Value = 145.875
$ws.Range("N91").Value = -3783
$ws.Range("H123").Value = 128078
$ws.Range("J123").Value = 128078
$ws.Range("L123").Value = 128078
$ws.Range("N123").Value = -137878
$ws.Range("H131").Value = 1292.9
$ws.Range("J131").Value = 1670
$ws.Range("L131").Value = 5010
$ws.Range("N131").Value = -15090
$ws.Range("H137").Value = 3789.9
$ws.Range("J137").Value = 6474.25
$ws.Range("L137").Value = 19422.75
$ws.Range("N137").Value = -24522.75

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6340.5454
$ws.Range("I74").Value = 4251.2856
$ws.Range("K74").Value = 4251.2856
$ws.Range("M74").Value = -3377.2856
$ws.Range("H77").Value = 6340.5454
$ws.Range("I77").Value = 4251.2856
$ws.Range("K77").Value = 21256.428
$ws.Range("M77").Value = -16888.428
$ws.Range("H122").Value = 1577.7273
$ws.Range("I122").Value = 1577.7273
$ws.Range("K122").Value = 4733.1819
$ws.Range("M122").Value = -2283.1819

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 15000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H64").Value = 441.5
$ws.Range("I64").Value = 315.5
$ws.Range("K64").Value = 315.5
$ws.Range("M64").Value = -90.5
$ws.Range("H67").Value = 441.5
$ws.Range("I67").Value = 315.5
$ws.Range("K67").Value = 315.5
$ws.Range("M67").Value = 464.5

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("H22").Value = 2241.6667
$ws.Range("J22").Value = 2999.5
$ws.Range("L22").Value = 2999.5
$ws.Range("N22").Value = -3699.5
$ws.Range("H31").Value = 8035.467
$ws.Range("I31").Value = 8035.467
$ws.Range("K31").Value = 8035.467
$ws.Range("M31").Value = -7740.467
$ws.Range("H34").Value = 8035.467
$ws.Range("I34").Value = 8035.467
$ws.Range("K34").Value = 8035.467
$ws.Range("M34").Value = -7833.467
$ws.Range("H38").Value = 39
$ws.Range("I38").Value = 38
$ws.Range("J38").Value = 42
$ws.Range("K38").Value = 38
$ws.Range("L38").Value = 42
$ws.Range("M38").Value = 339
$ws.Range("N38").Value = -796
$ws.Range("H46").Value = 39
$ws.Range("I46").Value = 38
$ws.Range("J46").Value = 42
$ws.Range("K46").Value = 38
$ws.Range("L46").Value = 42
$ws.Range("M46").Value = 173
$ws.Range("N46").Value = -464
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H134").Value = 2706.875
$ws.Range("I134").Value = 2098
$ws.Range("J134").Value = 6969
$ws.Range("K134").Value = 6294
$ws.Range("L134").Value = 20907
$ws.Range("M134").Value = -3759
$ws.Range("N134").Value = -25977
$ws.Range("H141").Value = 43430.938
$ws.Range("J141").Value = 43430.938
$ws.Range("L141").Value = 43430.938
$ws.Range("N141").Value = -53790.938

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 395.66666
$ws.Range("I8").Value = 395.66666
$ws.Range("K8").Value = 1186.99998
$ws.Range("M8").Value = -1047.99998
$ws.Range("H107").Value = 265.09525
$ws.Range("I107").Value = 429.8
$ws.Range("J107").Value = 213.625
$ws.Range("K107").Value = 1289.4
$ws.Range("L107").Value = 640.875
$ws.Range("M107").Value = 630.5999999999999
$ws.Range("N107").Value = -4480.875
$ws.Range("H140").Value = 3037.5
$ws.Range("I140").Value = 2153.5715
$ws.Range("K140").Value = 6460.7145
$ws.Range("M140").Value = -1280.7145

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2195.9
$ws.Range("I43").Value = 869.875
$ws.Range("K43").Value = 869.875
$ws.Range("M43").Value = -718.875
$ws.Range("H70").Value = 7627.75
$ws.Range("I70").Value = 6750.75
$ws.Range("K70").Value = 6750.75
$ws.Range("M70").Value = -6480.75
$ws.Range("H73").Value = 7627.75
$ws.Range("I73").Value = 6750.75
$ws.Range("K73").Value = 6750.75
$ws.Range("M73").Value = -5814.75
$ws.Range("H102").Value = 1189.619
$ws.Range("I102").Value = 567.64703
$ws.Range("J102").Value = 3833
$ws.Range("K102").Value = 567.64703
$ws.Range("L102").Value = 3833
$ws.Range("M102").Value = 1054.35297
$ws.Range("N102").Value = -7077
$ws.Range("H132").Value = 59815.19
$ws.Range("I132").Value = 76320.06
$ws.Range("K132").Value = 228960.18
$ws.Range("M132").Value = -226430.18

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1131
$ws.Range("J22").Value = 1497
$ws.Range("L22").Value = 1497
$ws.Range("N22").Value = -2087
$ws.Range("H27").Value = 1131
$ws.Range("J27").Value = 1497
$ws.Range("L27").Value = 1497
$ws.Range("N27").Value = -1711
$ws.Range("H46").Value = 6199.25
$ws.Range("J46").Value = 6999.2856
$ws.Range("L46").Value = 6999.2856
$ws.Range("N46").Value = -7375.2856

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8750
$ws.Range("I62").Value = 6166.6665
$ws.Range("J62").Value = 10300
$ws.Range("K62").Value = 6166.6665
$ws.Range("L62").Value = 10300
$ws.Range("M62").Value = -5542.6665
$ws.Range("N62").Value = -11548
$ws.Range("H65").Value = 8750
$ws.Range("I65").Value = 6166.6665
$ws.Range("J65").Value = 10300
$ws.Range("K65").Value = 30833.3325
$ws.Range("L65").Value = 51500
$ws.Range("M65").Value = -27713.3325
$ws.Range("N65").Value = -57740
